# "fixes inventory constraint bug"
# Updates the Part A / Part B revenue, material cost, inventory cost and
# profit figures on Sheet1 (the dependent "Percent Changed" formulas in
# column F recalculate automatically), re-applies the highlight style to
# the corrected "Part A" column (D6:D9), restores the selection to F9,
# and resizes the chart to its new (shorter) extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- corrected figures -----------------------------------------------
# revenue (row 6) - Part A unchanged, Part B corrected
$ws.Range("E6").Value = 391500

# MeterialCost (row 7)
$ws.Range("D7").Value = 240907
$ws.Range("E7").Value = 236796.29629629629

# InventoryCost (row 8)
$ws.Range("D8").Value = 56250
$ws.Range("E8").Value = 54425.000000000007

# Profit (row 9)
$ws.Range("D9").Value = 107843
$ws.Range("E9").Value = 100278.70370370371

# ---- re-style the "Part A" column (D6:D9) -----------------------------
# Apply formatting to D6 first, then copy that exact format down to
# D7:D9 so a single new cell style is produced (font color black, fill
# light green) instead of one per cell.
$d6 = $ws.Range("D6")
$d6.Font.Color = 0
$d6.Interior.Color = 14348258
$d6.Copy() | Out-Null
$ws.Range("D7:D9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- selection / view -------------------------------------------------
$ws.Range("F9").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1

# ---- resize the chart to its new (shorter) extent ---------------------
$co = $ws.ChartObjects().Item(1)
$co.Height = 202.2
